$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("registered_user")

# Rotate columns A (fname), B (lname), C (email) -> A (email), B (fname), C (lname)
# without disturbing the worksheet's <cols> column-width metadata (a whole-column
# Cut/Insert would otherwise stamp explicit default widths on columns A-C).
for ($r = 1; $r -le 21; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()
    $c = $ws.Cells.Item($r, 3).Value()

    $ws.Cells.Item($r, 1).Value = $c
    $ws.Cells.Item($r, 2).Value = $a
    $ws.Cells.Item($r, 3).Value = $b
}

# Add a new empty column F (after the existing E/bank_info column) with the
# same style/formatting as the rest of the text columns (A1:A21 uses style index 1).
$ws.Range("A1:A21").Copy() | Out-Null
$ws.Range("F1:F21").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("F1:F21").ClearContents() | Out-Null
$ws.Application.CutCopyMode = $false

# Update the selected cell as recorded in the saved view state.
$ws.Range("H8").Select() | Out-Null
